$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "29.540.18"
$ws.Range("E2").Value = "  +2.10%  "

# Row 3
$ws.Range("D3").Value = "1.988.49"
$ws.Range("E3").Value = "  +5.18%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9988"
$ws.Range("E4").Value = "  -0.30%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "328.29"
$ws.Range("E5").Value = "  +0.65%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9978"
$ws.Range("E6").Value = "  -0.41%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4650"
$ws.Range("E7").Value = "  +1.82%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3953"
$ws.Range("E8").Value = "  +1.46%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "46.11"
$ws.Range("E9").Value = "  -1.27%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07903"
$ws.Range("E10").Value = "  +0.71%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.005"
$ws.Range("E11").Value = "  +1.95%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "22.43"
$ws.Range("E12").Value = "  +2.66%  "

# Row 13
$ws.Range("D13").Value = "1.976.99"
$ws.Range("E13").Value = "  +3.27%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.221"
$ws.Range("E14").Value = "  +2.57%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.866"
$ws.Range("E15").Value = "  +3.25%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.07138"
$ws.Range("E16").Value = "  +2.72%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "88.88"
$ws.Range("E17").Value = "  +0.99%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.000"
$ws.Range("E18").Value = "  -0.26%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000009992"
$ws.Range("E19").Value = "  +0.17%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.15"
$ws.Range("E20").Value = "  +1.05%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.009"
$ws.Range("E21").Value = "  +0.59%  "

# Row 22
$ws.Range("D22").Value = "29.614.87"
$ws.Range("E22").Value = "  +2.35%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.559"
$ws.Range("E23").Value = "  +5.19%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.25"
$ws.Range("E24").Value = "  +2.58%  "

# Row 25
$ws.Range("D25").Value = "2.222.42"
$ws.Range("E25").Value = "  +3.94%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.120"
$ws.Range("E26").Value = "  +3.20%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "157.92"
$ws.Range("E27").Value = "  +1.39%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.62"
$ws.Range("E28").Value = "  +1.68%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.013"
$ws.Range("E29").Value = "  +1.94%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "120.25"
$ws.Range("E30").Value = "  +2.22%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.904"
$ws.Range("E31").Value = "  -1.21%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09410"
$ws.Range("E32").Value = "  +1.11%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.9003"
$ws.Range("E33").Value = "  -0.34%  "

# Row 34
$ws.Range("B34").Value = "ARBITRUM"
$ws.Range("C34").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.357"
$ws.Range("E34").Value = "  +2.62%  "

# Row 35
$ws.Range("B35").Value = "Filecoin"
$ws.Range("C35").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.277"
$ws.Range("E35").Value = "  -0.25%  "

# Row 36
$ws.Range("B36").Value = "PEPE"
$ws.Range("C36").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.000003782"
$ws.Range("E36").Value = "  +129.69%  "

# Row 37
$ws.Range("B37").Value = "HuobiToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.172"
$ws.Range("E37").Value = "  -2.82%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.05831"
$ws.Range("E38").Value = "  +1.07%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.182"
$ws.Range("E39").Value = "  -1.56%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.02131"
$ws.Range("E40").Value = "  +3.19%  "

# Row 41
$ws.Range("B41").Value = "FraxShare"
$ws.Range("C41").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "7.857"
$ws.Range("E41").Value = "  +2.41%  "

# Row 42
$ws.Range("B42").Value = "Frax"
$ws.Range("C42").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9978"
$ws.Range("E42").Value = "  -0.45%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.5780"
$ws.Range("E43").Value = "  +1.89%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.1826"
$ws.Range("E44").Value = "  +3.01%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "9.847"
$ws.Range("E45").Value = "  +1.14%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.21"
$ws.Range("E46").Value = "  +2.73%  "

# Row 47
$ws.Range("B47").Value = "RenderToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.233"
$ws.Range("E47").Value = "  -2.39%  "

# Row 48
$ws.Range("B48").Value = "Decentraland"
$ws.Range("C48").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.5390"
$ws.Range("E48").Value = "  +0.69%  "

# Row 49
$ws.Range("B49").Value = "MXToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.642"
$ws.Range("E49").Value = "  +5.22%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.07001"
$ws.Range("E50").Value = "  -0.53%  "

# Row 51
$ws.Range("B51").Value = "NEARProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.881"
$ws.Range("E51").Value = "  +2.05%  "
